$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '31.029.72'
$ws.Range('E2').Value = '  +0.79%  '

$ws.Range('D3').Value = '1.978.32'
$ws.Range('E3').Value = '  +0.51%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.005'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  +1.59%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '252.82'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.52%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.7313'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +2.30%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.006'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +1.70%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3365'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +3.89%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '27.49'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +6.58%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.07100'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +3.72%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.8248'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -0.80%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.08097'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +2.12%  '

$ws.Range('D13').Value = '1.984.64'
$ws.Range('E13').Value = '  +1.98%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.558'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +3.42%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '98.79'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -2.57%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '15.23'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +10.10%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '266.87'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -3.27%  '

$ws.Range('D18').Value = '31.042.20'
$ws.Range('E18').Value = '  +0.86%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.066'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +8.20%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.000008149'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +5.80%  '

$ws.Range('D21').Value = '2.246.29'
$ws.Range('E21').Value = '  +2.15%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '1.006'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +1.43%  '

$ws.Range('E23').Value = '  +1.68%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '7.055'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +6.22%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '9.914'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +4.27%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '161.94'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -1.48%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '19.63'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.74%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.344'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +8.30%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.1324'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +3.31%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.590'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +2.82%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.373'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +1.79%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.617'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +3.02%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.393'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +1.41%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.05277'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +4.88%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.271'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +5.77%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.7745'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +6.11%  '

$ws.Range('E37').Value = '  +3.74%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01997'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +1.14%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.876'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -1.33%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '83.28'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +7.40%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '6.735'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +2.91%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.4617'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +0.14%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.082'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +1.51%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.8517'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +1.50%  '

$ws.Range('B45').Value = 'PaxDollar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.007'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +1.74%  '

$ws.Range('B46').Value = 'Quant'
$ws.Range('C46').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '104.55'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +2.49%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '10.04'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +0.39%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '7.632'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +3.72%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.578'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +11.57%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '37.12'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +3.14%  '

$ws.Range('B51').Value = 'Decentraland'
$ws.Range('C51').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.4283'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +1.96%  '
